$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B2" = 1.02
    "C2" = 1.016472322111586
    "D2" = 1.022566044205062
    "E2" = 1.020082903372248
    "F2" = 1.014817610740538
    "I2" = 1.026629750270149
    "J2" = 1.021691936902136
    "K2" = 1.025400048454079
    "L2" = 1.022924240489359
    "M2" = 1.017674625892047
    "N2" = 1.011308899802676
    "B3" = 1.02
    "C3" = 1.017477873993032
    "D3" = 1.023294433916697
    "E3" = 1.021035135247171
    "F3" = 1.016466835624503
    "I3" = 1.026763881074105
    "J3" = 1.022332659818028
    "K3" = 1.025935345463548
    "L3" = 1.023682234416088
    "M3" = 1.019126537224389
    "N3" = 1.011523902568712
    "B4" = 1.02
    "C4" = 1.018128356564159
    "D4" = 1.023765331290641
    "E4" = 1.021651490862059
    "F4" = 1.01753364878944
    "I4" = 1.026849118350119
    "J4" = 1.022746527042803
    "K4" = 1.026280635116439
    "L4" = 1.02417229136844
    "M4" = 1.020065227862756
    "N4" = 1.011662695540496
    "B5" = 1.02
    "C5" = 1.018401778038899
    "D5" = 1.023963195951036
    "E5" = 1.021910654573872
    "F5" = 1.017982061317517
    "I5" = 1.026884579888298
    "J5" = 1.022920344186586
    "K5" = 1.026425535648854
    "L5" = 1.024378212258574
    "M5" = 1.020459669586501
    "N5" = 1.011720965636522
    "B6" = 1.02
    "C6" = 1.018447684306448
    "D6" = 1.023996412394007
    "E6" = 1.021954172117488
    "F6" = 1.018057347516927
    "I6" = 1.026890512206775
    "J6" = 1.022949518741679
    "K6" = 1.026449849887026
    "L6" = 1.024412781487792
    "M6" = 1.020525887483533
    "N6" = 1.011730744843457
    "B7" = 1.02
    "C7" = 1.018132010196823
    "D7" = 1.023767975563432
    "E7" = 1.021654953631876
    "F7" = 1.017539640788439
    "I7" = 1.026849593651607
    "J7" = 1.022748850273686
    "K7" = 1.026282572302959
    "L7" = 1.024175043282067
    "M7" = 1.020070499127174
    "N7" = 1.011663474456582
    "B8" = 1.02
    "C8" = 1.016812190223494
    "D8" = 1.022812293673144
    "E8" = 1.020404673873604
    "F8" = 1.015375049801406
    "I8" = 1.026675402174423
    "J8" = 1.021908622307984
    "K8" = 1.025581178958053
    "L8" = 1.023180494479121
    "M8" = 1.018165473980264
    "N8" = 1.011381629010824
    "B9" = 1.02
    "C9" = 1.014485103668999
    "D9" = 1.021125054074081
    "E9" = 1.018203016925073
    "F9" = 1.011557802582353
    "I9" = 1.026356554546336
    "J9" = 1.020422468778862
    "K9" = 1.024336931806873
    "L9" = 1.021424765091927
    "M9" = 1.014802231285352
    "N9" = 1.010882462930287
    "B10" = 1.02
    "C10" = 1.012932713573918
    "D10" = 1.019998078023773
    "E10" = 1.016736219957202
    "F10" = 1.009010539670433
    "I10" = 1.026135993124235
    "J10" = 1.019427921756769
    "K10" = 1.023501841050101
    "L10" = 1.02025208781723
    "M10" = 1.012555423550064
    "N10" = 1.010547984612515
    "B11" = 1.02
    "C11" = 1.012260260053303
    "D11" = 1.019509575322183
    "E11" = 1.016101303806831
    "F11" = 1.007906871740349
    "I11" = 1.026038591519975
    "J11" = 1.018996366884445
    "K11" = 1.023138907842106
    "L11" = 1.019743777373642
    "M11" = 1.011581344162813
    "N11" = 1.010402746352187
    "B12" = 1.02
    "C12" = 1.012010440701822
    "D12" = 1.019328046255483
    "E12" = 1.015865499430317
    "F12" = 1.007496808707299
    "I12" = 1.026002127198937
    "J12" = 1.018835930694771
    "K12" = 1.02300389758628
    "L12" = 1.01955488726339
    "M12" = 1.01121934052614
    "N12" = 1.010348737047762
    "B13" = 1.02
    "C13" = 1.012064029640546
    "D13" = 1.019366988361243
    "E13" = 1.015916078845573
    "F13" = 1.00758477380829
    "I13" = 1.026009961818038
    "J13" = 1.018870351037941
    "K13" = 1.023032866832798
    "L13" = 1.019595408511232
    "M13" = 1.01129700012422
    "N13" = 1.010360325006422
    "B14" = 1.02
    "C14" = 1.012239610735175
    "D14" = 1.01949457165269
    "E14" = 1.016081811500418
    "F14" = 1.007872978128622
    "I14" = 1.026035583182431
    "J14" = 1.018983107988832
    "K14" = 1.02312775194872
    "L14" = 1.019728165315018
    "M14" = 1.011551424690516
    "N14" = 1.010398283178938
    "B15" = 1.02
    "C15" = 1.012347786678431
    "D15" = 1.019573169668349
    "E15" = 1.016183929042386
    "F15" = 1.008050535275621
    "I15" = 1.026051331581025
    "J15" = 1.019052563022366
    "K15" = 1.023186187193805
    "L15" = 1.019809950407699
    "M15" = 1.01170815903827
    "N15" = 1.01042166232629
    "B16" = 1.02
    "C16" = 1.0129773365404
    "D16" = 1.020030487488665
    "E16" = 1.016778361747512
    "F16" = 1.009083771291024
    "I16" = 1.026142417377776
    "J16" = 1.019456543398079
    "K16" = 1.02352589960688
    "L16" = 1.020285811397201
    "M16" = 1.012620044146985
    "N16" = 1.010557615005945
    "B17" = 1.02
    "C17" = 1.013372166410732
    "D17" = 1.020317213030085
    "E17" = 1.017151291117499
    "F17" = 1.009731703212715
    "I17" = 1.026199045185474
    "J17" = 1.01970970552546
    "K17" = 1.023738634978553
    "L17" = 1.02058416301985
    "M17" = 1.013191719972083
    "N17" = 1.010642785433778
    "B18" = 1.02
    "C18" = 1.013602439157171
    "D18" = 1.020484405539602
    "E18" = 1.017368835463697
    "F18" = 1.010109565298104
    "I18" = 1.026231892190912
    "J18" = 1.019857283020871
    "K18" = 1.023862591293585
    "L18" = 1.020758135060509
    "M18" = 1.013525053830584
    "N18" = 1.010692424653677
    "B19" = 1.02
    "C19" = 1.013680952041254
    "D19" = 1.020541405446276
    "E19" = 1.01744301604202
    "F19" = 1.010238395628735
    "I19" = 1.026243061125202
    "J19" = 1.019907588285291
    "K19" = 1.02390483538739
    "L19" = 1.020817446335207
    "M19" = 1.013638692783656
    "N19" = 1.010709343693753
    "B20" = 1.02
    "C20" = 1.013329807480826
    "D20" = 1.020286455228858
    "E20" = 1.017111277163594
    "F20" = 1.009662193099778
    "I20" = 1.02619298849004
    "J20" = 1.019682552701342
    "K20" = 1.023715823793633
    "L20" = 1.020552158038042
    "M20" = 1.013130396511523
    "N20" = 1.010633651514146
    "B21" = 1.02
    "C21" = 1.012187907598965
    "D21" = 1.019457003691081
    "E21" = 1.016033006494277
    "F21" = 1.007788112289557
    "I21" = 1.026028046197552
    "J21" = 1.018949907679689
    "K21" = 1.023099816190466
    "L21" = 1.019689073987355
    "M21" = 1.011476508223645
    "N21" = 1.010387107138416
    "B22" = 1.02
    "C22" = 1.011469716714137
    "D22" = 1.018935046646382
    "E22" = 1.015355238230583
    "F22" = 1.006609153235295
    "I22" = 1.02592269133385
    "J22" = 1.018488468738044
    "K22" = 1.022711346364103
    "L22" = 1.019145949848278
    "M22" = 1.010435556360838
    "N22" = 1.010231739858003
    "B23" = 1.02
    "C23" = 1.011850465842513
    "D23" = 1.0192117884462
    "E23" = 1.015714518903987
    "F23" = 1.007234206212462
    "I23" = 1.0259786982628
    "J23" = 1.018733161987703
    "K23" = 1.022917391783946
    "L23" = 1.019433914962615
    "M23" = 1.01098749011452
    "N23" = 1.010314136718288
    "B24" = 1.02
    "C24" = 1.013348947734252
    "D24" = 1.02030035350851
    "E24" = 1.017129357681929
    "F24" = 1.009693601928644
    "I24" = 1.026195725815482
    "J24" = 1.01969482216443
    "K24" = 1.023726131583435
    "L24" = 1.020566619870235
    "M24" = 1.013158106297246
    "N24" = 1.01063777886045
    "B25" = 1.02
    "C25" = 1.015086882906321
    "D25" = 1.021561625230809
    "E25" = 1.018772023951865
    "F25" = 1.012545048641716
    "I25" = 1.026440394029039
    "J25" = 1.020807338259856
    "K25" = 1.024659584455405
    "L25" = 1.021879046832401
    "M25" = 1.015672501062421
    "N25" = 1.011011808549178
}

foreach ($cell in $values.Keys) {
    $ws.Range($cell).Value = $values[$cell]
}
